# Refresh the cryptos list: updated Price (D) / Volume(1h) (E) figures, plus
# two pairs of rows (33/34, 41/42, 45/46) whose coins were re-ranked and swapped
# places with their neighbour, and row 51 whose coin was replaced entirely.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D stores prices as text (e.g. "98.905.08"). Several of the new prices
# (e.g. "262.26") also look like plain numbers, and a plain Value assignment
# would make Excel silently convert the cell to a Number (losing the exact
# text, e.g. "1.40" -> 1.4). Force those cells to Text first, write the value,
# then restore the Normal style so no stray number formatting is left behind.
function Set-TextValue($cell, [string]$value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

$ws.Range('D2').Value = '98.905.08'
$ws.Range('E2').Value = '  +2.63%  '

$ws.Range('D3').Value = '3.411.83'
$ws.Range('E3').Value = '  +9.25%  '

$ws.Range('E4').Value = '  -0.10%  '

Set-TextValue $ws.Range('D5') '262.26'
$ws.Range('E5').Value = '  +9.46%  '

Set-TextValue $ws.Range('D6') '636.71'
$ws.Range('E6').Value = '  +4.86%  '

Set-TextValue $ws.Range('D7') '1.40'
$ws.Range('E7').Value = '  +27.21%  '

$ws.Range('E8').Value = '  +2.83%  '

$ws.Range('E9').Value = '  -0.09%  '

Set-TextValue $ws.Range('D10') '0.891'
$ws.Range('E10').Value = '  +13.12%  '

$ws.Range('D11').Value = '3.412.91'
$ws.Range('E11').Value = '  +9.58%  '

Set-TextValue $ws.Range('D12') '0.201'
$ws.Range('E12').Value = '  +2.10%  '

$ws.Range('D13').Value = '98.546.88'
$ws.Range('E13').Value = '  +2.77%  '

Set-TextValue $ws.Range('D14') '36.52'
$ws.Range('E14').Value = '  +7.94%  '

Set-TextValue $ws.Range('D15') '0.0000251'
$ws.Range('E15').Value = '  +5.02%  '

$ws.Range('D16').Value = '4.024.51'
$ws.Range('E16').Value = '  +9.26%  '

Set-TextValue $ws.Range('D17') '5.58'
$ws.Range('E17').Value = '  +4.67%  '

$ws.Range('D18').Value = '3.381.98'
$ws.Range('E18').Value = '  +9.32%  '

Set-TextValue $ws.Range('D19') '3.64'
$ws.Range('E19').Value = '  +2.65%  '

Set-TextValue $ws.Range('D20') '15.28'
$ws.Range('E20').Value = '  +6.01%  '

Set-TextValue $ws.Range('D21') '493.82'
$ws.Range('E21').Value = '  +5.43%  '

Set-TextValue $ws.Range('D22') '6.23'
$ws.Range('E22').Value = '  +9.89%  '

Set-TextValue $ws.Range('D23') '0.0000217'
$ws.Range('E23').Value = '  +13.55%  '

Set-TextValue $ws.Range('D24') '9.48'
$ws.Range('E24').Value = '  +8.38%  '

Set-TextValue $ws.Range('D25') '5.77'
$ws.Range('E25').Value = '  +4.73%  '

Set-TextValue $ws.Range('D26') '88.98'
$ws.Range('E26').Value = '  +4.46%  '

Set-TextValue $ws.Range('D27') '12.17'
$ws.Range('E27').Value = '  +4.79%  '

$ws.Range('D28').Value = '3.537.86'
$ws.Range('E28').Value = '  +8.39%  '

$ws.Range('E29').Value = '  +21.12%  '

Set-TextValue $ws.Range('D30') '0.999'
$ws.Range('E30').Value = '  -0.06%  '

Set-TextValue $ws.Range('D31') '0.195'
$ws.Range('E31').Value = '  +11.07%  '

Set-TextValue $ws.Range('D32') '0.132'
$ws.Range('E32').Value = '  +5.81%  '

$ws.Range('B33').Value = 'InternetComputer(DFINITY)'
$ws.Range('C33').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue $ws.Range('D33') '9.68'
$ws.Range('E33').Value = '  +7.12%  '

$ws.Range('B34').Value = 'Binance-PegBSC-USD'
$ws.Range('C34').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
Set-TextValue $ws.Range('D34') '1.00'
$ws.Range('E34').Value = '  +0.23%  '

Set-TextValue $ws.Range('D35') '28.27'
$ws.Range('E35').Value = '  +8.15%  '

Set-TextValue $ws.Range('D36') '7.46'
$ws.Range('E36').Value = '  +1.54%  '

Set-TextValue $ws.Range('D37') '0.152'
$ws.Range('E37').Value = '  +1.05%  '

Set-TextValue $ws.Range('D38') '2.01'
$ws.Range('E38').Value = '  +7.53%  '

Set-TextValue $ws.Range('D39') '513.18'
$ws.Range('E39').Value = '  +5.26%  '

Set-TextValue $ws.Range('D40') '0.476'
$ws.Range('E40').Value = '  +8.90%  '

$ws.Range('B41').Value = 'WhiteBITCoin'
$ws.Range('C41').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
Set-TextValue $ws.Range('D41') '24.85'
$ws.Range('E41').Value = '  +2.92%  '

$ws.Range('B42').Value = 'MantraDAO'
$ws.Range('C42').Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
Set-TextValue $ws.Range('D42') '3.83'
$ws.Range('E42').Value = '  +5.81%  '

Set-TextValue $ws.Range('D43') '1.29'
$ws.Range('E43').Value = '  +5.14%  '

Set-TextValue $ws.Range('D44') '3.42'
$ws.Range('E44').Value = '  +8.12%  '

$ws.Range('B45').Value = 'ARBITRUM'
$ws.Range('C45').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue $ws.Range('D45') '0.796'
$ws.Range('E45').Value = '  +14.29%  '

$ws.Range('B46').Value = 'USDe'
$ws.Range('C46').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
Set-TextValue $ws.Range('D46') '1.00'
$ws.Range('E46').Value = '  +0.03%  '

Set-TextValue $ws.Range('D47') '161.77'
$ws.Range('E47').Value = '  +0.00%  '

Set-TextValue $ws.Range('D48') '1.97'
$ws.Range('E48').Value = '  +4.66%  '

Set-TextValue $ws.Range('D49') '4.70'
$ws.Range('E49').Value = '  +8.64%  '

Set-TextValue $ws.Range('D50') '46.71'
$ws.Range('E50').Value = '  +6.60%  '

$ws.Range('B51').Value = 'Mantle'
$ws.Range('C51').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextValue $ws.Range('D51') '0.824'
$ws.Range('E51').Value = '  +12.42%  '
